$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "Sheet1"

# --- Header row (row 1) ---
$ws.Cells.Item(1,1).Value2 = "Product"
$ws.Cells.Item(1,2).Value2 = "Price"
$ws.Cells.Item(1,3).Value2 = "Data"

# --- Data rows (2-9): product name / price / scrape timestamp ---
$ws.Cells.Item(2,1).Value2 = "Cafea Organica House Roast, Exhale, boabe"
$ws.Cells.Item(2,2).Value2 = "155,00 Lei"
$ws.Cells.Item(2,3).Value2 = 45504.042433495371

$ws.Cells.Item(3,1).Value2 = "Klorane Urzica Sampon uscat 150 ml"
$ws.Cells.Item(3,2).Value2 = "35,05 Lei"
$ws.Cells.Item(3,3).Value2 = 45504.042444629631

$ws.Cells.Item(4,1).Value2 = "Crema pentru ochi Elmiplant Hyaluronic Gold, Femei, 15 ml"
$ws.Cells.Item(4,2).Value2 = "31,40 Lei"
$ws.Cells.Item(4,3).Value2 = 45504.042452569447

$ws.Cells.Item(5,1).Value2 = "Bautura de Ovaz Standard Minor Figures bax 6L"
$ws.Cells.Item(5,2).Value2 = "108,00 Lei"
$ws.Cells.Item(5,3).Value2 = 45504.042459872682

$ws.Cells.Item(6,1).Value2 = "Cafea Organica House Roast, Exhale, boabe"
$ws.Cells.Item(6,2).Value2 = "155,00 Lei"
$ws.Cells.Item(6,3).Value2 = 45504.044278368063
$ws.Cells.Item(6,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

$ws.Cells.Item(7,1).Value2 = "Klorane Urzica Sampon uscat 150 ml"
$ws.Cells.Item(7,2).Value2 = "35,05 Lei"
$ws.Cells.Item(7,3).Value2 = 45504.044287997684
$ws.Cells.Item(7,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

$ws.Cells.Item(8,1).Value2 = "Crema pentru ochi Elmiplant Hyaluronic Gold, Femei, 15 ml"
$ws.Cells.Item(8,2).Value2 = "31,40 Lei"
$ws.Cells.Item(8,3).Value2 = 45504.044296539352
$ws.Cells.Item(8,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

$ws.Cells.Item(9,1).Value2 = "Bautura de Ovaz Standard Minor Figures bax 6L"
$ws.Cells.Item(9,2).Value2 = "108,00 Lei"
$ws.Cells.Item(9,3).Value2 = 45504.044303949639
$ws.Cells.Item(9,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

# --- Header style: bold, centered (horizontal + vertical) ---
$header = $ws.Range("A1:C1")
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.Font.Bold = $true

# --- Drop the now-unused column formatting (old D/E/F widths) ---
$ws.Range("D1:F1").EntireColumn.Delete()

# --- Column widths for A, B, C ---
$ws.Columns.Item(1).ColumnWidth = 39
$ws.Columns.Item(2).ColumnWidth = 29.333333333333332
$ws.Columns.Item(3).ColumnWidth = 31.333333333333332

# --- View: reset zoom, move selection ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("F6").Select()
